$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = 2
$ws.Range("D3").Value = 5
$ws.Range("E3").Value = 0.04

$ws.Range("D4").Value = 1

$ws.Range("C5").Value = 5
$ws.Range("E5").Value = 0.15

$ws.Range("D6").Value = 6

$ws.Range("C7").Value = 0
$ws.Range("D7").Value = 5
$ws.Range("E7").Value = 0

$ws.Range("C9").Value = 39
$ws.Range("D9").Value = 37
$ws.Range("E9").Value = 0.06

$ws.Range("C10").Value = 0
$ws.Range("E10").Value = 1

$ws.Range("C11").Value = 2
$ws.Range("E11").Value = 0.27

$ws.Range("D12").Value = 9
$ws.Range("E12").Value = 0

$ws.Range("C13").Value = 88

$ws.Range("C14").Value = 5
$ws.Range("E14").Value = 0.01

$ws.Range("C17").Value = 11
$ws.Range("D17").Value = 20
$ws.Range("E17").Value = 0

$ws.Range("D19").Value = 3
$ws.Range("E19").Value = 0.09

$ws.Range("C21").Value = 1
$ws.Range("D21").Value = 1
$ws.Range("E21").Value = 0.37

$ws.Range("D22").Value = 0
$ws.Range("E22").Value = 1

$ws.Range("C24").Value = 5
$ws.Range("D24").Value = 6
$ws.Range("E24").Value = 0.15

$ws.Range("C25").Value = 1
$ws.Range("E25").Value = 0.37

$ws.Range("D28").Value = 1
$ws.Range("E28").Value = 0.27

$ws.Range("D29").Value = 0
$ws.Range("E29").Value = 1

$ws.Range("D30").Value = 1
$ws.Range("E30").Value = 0.37

$ws.Range("C32").Value = 5
$ws.Range("D32").Value = 5
$ws.Range("E32").Value = 0.18

$ws.Range("D33").Value = 1
$ws.Range("E33").Value = 0

$ws.Range("C34").Value = 11
$ws.Range("D34").Value = 8
$ws.Range("E34").Value = 0.09

